# CIV-17609 - updated GA documents to display main claim number
#
# The template labelled the case-number field as "Claim number:" in the
# first order block even though the merge field itself is `caseNumber`
# and the later block in the same document already reads
# "Case number: <<caseNumber>>". Correct the mislabelled heading so both
# instances read "Case number".

$d = $word.ActiveDocument

# Target the unique phrase so we don't touch "Claimant" elsewhere in the
# document.
$found = $d.Content.Find.Execute(
    "Claim number: <<",   # FindText
    $false,                # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                # Format
    "Case number: <<",    # ReplaceWith
    2                       # Replace (wdReplaceAll)
)

Write-Output ("Claim number -> Case number replaced: " + $found)
